$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 'LinkedIn Video - Awareness'
$ws.Range("B2").Value2 = 'LinkedIn Video Ad 1'
$ws.Range("C2").Value2 = 'ACTIVE'
$ws.Range("D2").Value2 = 'https://video.talentflow.com/enterprise/demo-1.mp4'
$ws.Range("E2").Value2 = 'Discover how leading enterprises are revolutionizing their hiring process with AI-powered recruitment. See real results from Fortune 500 companies.'
$ws.Range("F2").Value2 = 'https://talentflow.com/enterprise'
$ws.Range("G2").Value2 = 'LEARN_MORE'
$ws.Range("H2").Value2 = 75
$ws.Range("I2").Value2 = 6000
$ws.Range("J2").Value2 = 'https://cdn.talentflow.com/enterprise/thumbnail-1.jpg'

# Row 3
$ws.Range("A3").Value2 = 'LinkedIn Video - Consideration'
$ws.Range("B3").Value2 = 'LinkedIn Video Ad 2'
$ws.Range("C3").Value2 = 'ACTIVE'
$ws.Range("D3").Value2 = 'https://video.securecloud.com/business/demo-2.mp4'
$ws.Range("E3").Value2 = 'See why 10,000+ businesses trust our enterprise cloud security platform. SOC 2 certified, GDPR compliant, 99.99% uptime SLA.'
$ws.Range("F3").Value2 = 'https://securecloud.com/business'
$ws.Range("G3").Value2 = 'CONTACT_US'
$ws.Range("H3").Value2 = 50
$ws.Range("I3").Value2 = 6000
$ws.Range("J3").Value2 = 'https://cdn.securecloud.com/business/thumbnail-2.jpg'

# Row 4
$ws.Range("A4").Value2 = 'LinkedIn Video - Conversion'
$ws.Range("B4").Value2 = 'LinkedIn Video Ad 3'
$ws.Range("C4").Value2 = 'ACTIVE'
$ws.Range("D4").Value2 = 'https://video.datasync.com/platform/demo-3.mp4'
$ws.Range("E4").Value2 = 'Learn how our data integration platform helps companies connect 500+ applications seamlessly. No-code solution, real-time sync.'
$ws.Range("F4").Value2 = 'https://datasync.com/platform'
$ws.Range("G4").Value2 = 'DOWNLOAD'
$ws.Range("H4").Value2 = 200
$ws.Range("I4").Value2 = 2250
$ws.Range("J4").Value2 = 'https://cdn.datasync.com/platform/thumbnail-3.jpg'

# Row 5
$ws.Range("A5").Value2 = 'LinkedIn Video - Retargeting'
$ws.Range("B5").Value2 = 'LinkedIn Video Ad 4'
$ws.Range("C5").Value2 = 'ACTIVE'
$ws.Range("D5").Value2 = 'https://video.salesai.com/demo/demo-4.mp4'
$ws.Range("E5").Value2 = 'Watch how AI is transforming B2B sales. Our platform helped customers increase revenue by 45% on average in the first year.'
$ws.Range("F5").Value2 = 'https://salesai.com/demo'
$ws.Range("G5").Value2 = 'REQUEST_DEMO'
$ws.Range("H5").Value2 = 200
$ws.Range("I5").Value2 = 3000
$ws.Range("J5").Value2 = 'https://cdn.salesai.com/demo/thumbnail-4.jpg'

# Row 6
$ws.Range("A6").Value2 = 'LinkedIn Video - Awareness'
$ws.Range("B6").Value2 = 'LinkedIn Video Ad 5'
$ws.Range("C6").Value2 = 'ACTIVE'
$ws.Range("D6").Value2 = 'https://video.teamhub.com/teams/demo-5.mp4'
$ws.Range("E6").Value2 = 'See how global teams collaborate in real-time with our enterprise workspace platform. Used by 5,000+ distributed teams worldwide.'
$ws.Range("F6").Value2 = 'https://teamhub.com/teams'
$ws.Range("G6").Value2 = 'CONTACT_US'
$ws.Range("H6").Value2 = 200
$ws.Range("I6").Value2 = 4500
$ws.Range("J6").Value2 = 'https://cdn.teamhub.com/teams/thumbnail-5.jpg'

# Row 7
$ws.Range("A7").Value2 = 'LinkedIn Video - Consideration'
$ws.Range("B7").Value2 = 'LinkedIn Video Ad 6'
$ws.Range("C7").Value2 = 'ACTIVE'
$ws.Range("D7").Value2 = 'https://video.financeos.com/cfo/demo-6.mp4'
$ws.Range("E7").Value2 = 'Discover why CFOs at growth companies choose our financial operations platform. Automated reporting, forecasting, and compliance.'
$ws.Range("F7").Value2 = 'https://financeos.com/cfo'
$ws.Range("G7").Value2 = 'REGISTER'
$ws.Range("H7").Value2 = 100
$ws.Range("I7").Value2 = 1500
$ws.Range("J7").Value2 = 'https://cdn.financeos.com/cfo/thumbnail-6.jpg'

# Row 8
$ws.Range("A8").Value2 = 'LinkedIn Video - Conversion'
$ws.Range("B8").Value2 = 'LinkedIn Video Ad 7'
$ws.Range("C8").Value2 = 'ACTIVE'
$ws.Range("D8").Value2 = 'https://video.leadgenpro.com/marketing/demo-7.mp4'
$ws.Range("E8").Value2 = 'Learn how marketing teams generate 10x more qualified leads with our B2B lead generation platform. ROI guaranteed.'
$ws.Range("F8").Value2 = 'https://leadgenpro.com/marketing'
$ws.Range("G8").Value2 = 'REQUEST_DEMO'
$ws.Range("H8").Value2 = 75
$ws.Range("I8").Value2 = 1500
$ws.Range("J8").Value2 = 'https://cdn.leadgenpro.com/marketing/thumbnail-7.jpg'

# Row 9
$ws.Range("A9").Value2 = 'LinkedIn Video - Retargeting'
$ws.Range("B9").Value2 = 'LinkedIn Video Ad 8'
$ws.Range("C9").Value2 = 'ACTIVE'
$ws.Range("D9").Value2 = 'https://video.cybershield.com/security/demo-8.mp4'
$ws.Range("E9").Value2 = 'Watch how enterprises protect their data with our advanced cybersecurity platform. Real-time threat detection and response.'
$ws.Range("F9").Value2 = 'https://cybershield.com/security'
$ws.Range("G9").Value2 = 'CONTACT_US'
$ws.Range("H9").Value2 = 50
$ws.Range("I9").Value2 = 3000
$ws.Range("J9").Value2 = 'https://cdn.cybershield.com/security/thumbnail-8.jpg'

# Row 10
$ws.Range("A10").Value2 = 'LinkedIn Video - Awareness'
$ws.Range("B10").Value2 = 'LinkedIn Video Ad 9'
$ws.Range("C10").Value2 = 'ACTIVE'
$ws.Range("D10").Value2 = 'https://video.talentflow.com/enterprise/demo-9.mp4'
$ws.Range("E10").Value2 = 'Discover how leading enterprises are revolutionizing their hiring process with AI-powered recruitment. See real results from Fortune 500 companies.'
$ws.Range("F10").Value2 = 'https://talentflow.com/enterprise'
$ws.Range("G10").Value2 = 'REGISTER'
$ws.Range("H10").Value2 = 50
$ws.Range("I10").Value2 = 1500
$ws.Range("J10").Value2 = 'https://cdn.talentflow.com/enterprise/thumbnail-9.jpg'

# Row 11
$ws.Range("A11").Value2 = 'LinkedIn Video - Consideration'
$ws.Range("B11").Value2 = 'LinkedIn Video Ad 10'
$ws.Range("C11").Value2 = 'ACTIVE'
$ws.Range("D11").Value2 = 'https://video.securecloud.com/business/demo-10.mp4'
$ws.Range("E11").Value2 = 'See why 10,000+ businesses trust our enterprise cloud security platform. SOC 2 certified, GDPR compliant, 99.99% uptime SLA.'
$ws.Range("F11").Value2 = 'https://securecloud.com/business'
$ws.Range("G11").Value2 = 'DOWNLOAD'
$ws.Range("H11").Value2 = 50
$ws.Range("I11").Value2 = 3000
$ws.Range("J11").Value2 = 'https://cdn.securecloud.com/business/thumbnail-10.jpg'

# Row 12
$ws.Range("A12").Value2 = 'LinkedIn Video - Conversion'
$ws.Range("B12").Value2 = 'LinkedIn Video Ad 11'
$ws.Range("C12").Value2 = 'ACTIVE'
$ws.Range("D12").Value2 = 'https://video.datasync.com/platform/demo-11.mp4'
$ws.Range("E12").Value2 = 'Learn how our data integration platform helps companies connect 500+ applications seamlessly. No-code solution, real-time sync.'
$ws.Range("F12").Value2 = 'https://datasync.com/platform'
$ws.Range("G12").Value2 = 'LEARN_MORE'
$ws.Range("H12").Value2 = 200
$ws.Range("I12").Value2 = 3000
$ws.Range("J12").Value2 = 'https://cdn.datasync.com/platform/thumbnail-11.jpg'

# Row 13
$ws.Range("A13").Value2 = 'LinkedIn Video - Retargeting'
$ws.Range("B13").Value2 = 'LinkedIn Video Ad 12'
$ws.Range("C13").Value2 = 'ACTIVE'
$ws.Range("D13").Value2 = 'https://video.salesai.com/demo/demo-12.mp4'
$ws.Range("E13").Value2 = 'Watch how AI is transforming B2B sales. Our platform helped customers increase revenue by 45% on average in the first year.'
$ws.Range("F13").Value2 = 'https://salesai.com/demo'
$ws.Range("G13").Value2 = 'LEARN_MORE'
$ws.Range("H13").Value2 = 100
$ws.Range("I13").Value2 = 4500
$ws.Range("J13").Value2 = 'https://cdn.salesai.com/demo/thumbnail-12.jpg'

# Row 14
$ws.Range("A14").Value2 = 'LinkedIn Video - Awareness'
$ws.Range("B14").Value2 = 'LinkedIn Video Ad 13'
$ws.Range("C14").Value2 = 'ACTIVE'
$ws.Range("D14").Value2 = 'https://video.teamhub.com/teams/demo-13.mp4'
$ws.Range("E14").Value2 = 'See how global teams collaborate in real-time with our enterprise workspace platform. Used by 5,000+ distributed teams worldwide.'
$ws.Range("F14").Value2 = 'https://teamhub.com/teams'
$ws.Range("G14").Value2 = 'CONTACT_US'
$ws.Range("H14").Value2 = 75
$ws.Range("I14").Value2 = 2250
$ws.Range("J14").Value2 = 'https://cdn.teamhub.com/teams/thumbnail-13.jpg'

# Row 15
$ws.Range("A15").Value2 = 'LinkedIn Video - Consideration'
$ws.Range("B15").Value2 = 'LinkedIn Video Ad 14'
$ws.Range("C15").Value2 = 'ACTIVE'
$ws.Range("D15").Value2 = 'https://video.financeos.com/cfo/demo-14.mp4'
$ws.Range("E15").Value2 = 'Discover why CFOs at growth companies choose our financial operations platform. Automated reporting, forecasting, and compliance.'
$ws.Range("F15").Value2 = 'https://financeos.com/cfo'
$ws.Range("G15").Value2 = 'LEARN_MORE'
$ws.Range("H15").Value2 = 50
$ws.Range("I15").Value2 = 3000
$ws.Range("J15").Value2 = 'https://cdn.financeos.com/cfo/thumbnail-14.jpg'

# Row 16
$ws.Range("A16").Value2 = 'LinkedIn Video - Conversion'
$ws.Range("B16").Value2 = 'LinkedIn Video Ad 15'
$ws.Range("C16").Value2 = 'ACTIVE'
$ws.Range("D16").Value2 = 'https://video.leadgenpro.com/marketing/demo-15.mp4'
$ws.Range("E16").Value2 = 'Learn how marketing teams generate 10x more qualified leads with our B2B lead generation platform. ROI guaranteed.'
$ws.Range("F16").Value2 = 'https://leadgenpro.com/marketing'
$ws.Range("G16").Value2 = 'CONTACT_US'
$ws.Range("H16").Value2 = 50
$ws.Range("I16").Value2 = 1500
$ws.Range("J16").Value2 = 'https://cdn.leadgenpro.com/marketing/thumbnail-15.jpg'

# Row 17
$ws.Range("A17").Value2 = 'LinkedIn Video - Retargeting'
$ws.Range("B17").Value2 = 'LinkedIn Video Ad 16'
$ws.Range("C17").Value2 = 'ACTIVE'
$ws.Range("D17").Value2 = 'https://video.cybershield.com/security/demo-16.mp4'
$ws.Range("E17").Value2 = 'Watch how enterprises protect their data with our advanced cybersecurity platform. Real-time threat detection and response.'
$ws.Range("F17").Value2 = 'https://cybershield.com/security'
$ws.Range("G17").Value2 = 'DOWNLOAD'
$ws.Range("H17").Value2 = 75
$ws.Range("I17").Value2 = 3000
$ws.Range("J17").Value2 = 'https://cdn.cybershield.com/security/thumbnail-16.jpg'

# Row 18
$ws.Range("A18").Value2 = 'LinkedIn Video - Thought Leadership'
$ws.Range("B18").Value2 = 'LinkedIn Video Ad 17'
$ws.Range("C18").Value2 = 'ACTIVE'
$ws.Range("D18").Value2 = 'https://video.analyticsplatform.com/webinar-series.mp4'
$ws.Range("E18").Value2 = 'In today''s rapidly evolving business landscape, data-driven decision making has become absolutely critical for organizations seeking competitive advantage. Our comprehensive analytics platform empowers business leaders to harness the full potential of their data through advanced visualization, predictive modeling, and automated insights generation. Join thousands of enterprises who have already transformed their operations with our award-winning solution.'
$ws.Range("F18").Value2 = 'https://analyticsplatform.com/webinar'
$ws.Range("G18").Value2 = 'REGISTER'
$ws.Range("H18").Value2 = 100
$ws.Range("I18").Value2 = 3000
$ws.Range("J18").Value2 = 'https://cdn.analyticsplatform.com/thumbnail-webinar.jpg'

# Row 19
$ws.Range("A19").Value2 = 'LinkedIn Video - Product Launch'
$ws.Range("B19").Value2 = 'LinkedIn Video Ad 18'
$ws.Range("C19").Value2 = 'ACTIVE'
$ws.Range("D19").Value2 = 'https://video.hrtech.com/product-launch.mp4'
$ws.Range("E19").Value2 = 'Transform your HR operations with our all-in-one people management platform. Recruitment, onboarding, performance, payroll.'
$ws.Range("F19").Value2 = 'https://hrtech.com/platform'
$ws.Range("G19").Value2 = 'REQUEST_DEMO'
$ws.Range("H19").Value2 = 5
$ws.Range("I19").Value2 = 150
$ws.Range("J19").Value2 = 'https://cdn.hrtech.com/thumbnail-platform.jpg'

# Row 20
$ws.Range("A20").Value2 = 'LinkedIn Video - Case Study'
$ws.Range("B20").Value2 = 'LinkedIn Video Ad 19'
$ws.Range("C20").Value2 = 'ACTIVE'
$ws.Range("D20").Value2 = $null
$ws.Range("E20").Value2 = 'See how Fortune 500 companies achieve 99.9% customer satisfaction with our enterprise support platform. Real results, real ROI.'
$ws.Range("F20").Value2 = 'https://supportpro.com/case-studies'
$ws.Range("G20").Value2 = 'LEARN_MORE'
$ws.Range("H20").Value2 = 75
$ws.Range("I20").Value2 = 2250
$ws.Range("J20").Value2 = 'https://cdn.supportpro.com/thumbnail-case-study.jpg'

# Row 21
$ws.Range("A21").Value2 = 'LinkedIn Video - Brand Awareness'
$ws.Range("B21").Value2 = 'LinkedIn Video Ad 20'
$ws.Range("C21").Value2 = 'active'
$ws.Range("D21").Value2 = 'https://video.cloudinfra.com/infrastructure-demo.mp4'
$ws.Range("E21").Value2 = 'Discover why leading enterprises choose our cloud infrastructure. Scalable, secure, and cost-effective solutions for growing businesses.'
$ws.Range("F21").Value2 = 'https://cloudinfra.com/enterprise'
$ws.Range("G21").Value2 = 'CONTACT_US'
$ws.Range("H21").Value2 = 100
$ws.Range("I21").Value2 = 3000
$ws.Range("J21").Value2 = 'https://cdn.cloudinfra.com/thumbnail-infra.jpg'
